$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.8544378796644995, 0.8568587807719603),
    @(0.6463411907890931, 0.6508037365575856),
    @(0.616621464201531, 0.6196006649465076),
    @(0.8194018943245311, 0.8230693242384203),
    @(0.0805354307539307, 0.08849797441731577),
    @(0.96925716638757, 0.9640248577244662),
    @(1.298841514114648, 1.285869942372226),
    @(1.753943295583375, 1.744144833335414),
    @(0.6003913602122007, 0.603435823035954),
    @(1.310652557793934, 1.307704649210073),
    @(1.750731214159401, 1.738927247565891),
    @(1.733198432686658, 1.721132472644895),
    @(1.514918593871193, 1.506946016148416),
    @(0.9124105864632828, 0.9117524573638238),
    @(0.9652833326741588, 0.9594030572056931),
    @(0.6519566776739463, 0.6558769406337663),
    @(0.6436283712855808, 0.6450455588157792),
    @(1.597932677876217, 1.584705878787109),
    @(0.9069490511470674, 0.912686059895371),
    @(0.7101530737634343, 0.7138587298365591),
    @(0.1413705257254187, 0.1501976737139035),
    @(1.4913702792392, 1.474531607867737),
    @(0.681790817376478, 0.6852419275788471),
    @(0.8287310113894724, 0.8316630937330675),
    @(1.287761039043411, 1.281951533063581),
    @(0.3093444057557265, 0.3130104681774866),
    @(1.282200554881656, 1.28780669613771),
    @(0.7520441022150447, 0.7531660135440238),
    @(1.163541580838492, 1.157771411623605),
    @(0.6358393769419587, 0.6392176430679251),
    @(0.6906127335161945, 0.6900792722450098),
    @(0.7695585959571494, 0.7643755423797318),
    @(0.0468689095057648, 0.06971745368735227),
    @(1.0164714988704, 1.011373846772858),
    @(0.3763628084380087, 0.3839347535113792),
    @(1.30186507690716, 1.296988677021973),
    @(0.398555750789446, 0.405733442907917),
    @(0.7199799374896158, 0.7258200266249029),
    @(0.9171691541809092, 0.9143059459606744),
    @(0.5582211534847016, 0.5564744892975657),
    @(0.7880253087745966, 0.7937830085106274),
    @(0.8577585081917437, 0.8567376221429097)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $values[$i][0]
    $ws.Cells.Item($r, 2).Value = $values[$i][1]
}
